{"js": "// Fixed Search Box and updates Resumen Ejecutivo\n//\n// 1) Merge the split runs around `\u201d o \u201c` in the \"Edit Product\" / \"Delete\n//    Product\" bullet so the quote + \" o \" text lives in a single run.\n// 2) Merge the split runs around the trailing `\u201d.` in that same bullet.\n// 3) Remove the empty paragraph and the \"B\u00fasqueda de productos\u2026\" bullet\n//    at the end of the document (the search-box filter bullet was\n//    dropped from the write-up).\n\nconst body = context.document.body;\n\n// --- 1 & 2: collapse the cosmetic run-splits around the smart quotes ---\n// Scope the search to the specific paragraph that contains both quoted\n// phrases (\"Edit Product\" o \"Delete Product\") so we don't touch the many\n// other smart-quote runs elsewhere in the document.\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Edit Product\") !== -1 && t.indexOf(\"Delete Product\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // \u201d o \u201c  ->  single run (was split across two runs)\n  const quoteOr = targetParagraph.search(\"\u201d o \u201c\", { matchCase: true, matchWildcards: false });\n  quoteOr.load(\"items\");\n  await context.sync();\n  if (quoteOr.items.length > 0) {\n    quoteOr.items[0].insertText(\"\u201d o \u201c\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // \u201d.  ->  single run (was split across two runs)\n  const quotePeriod = targetParagraph.search(\"\u201d.\", { matchCase: true, matchWildcards: false });\n  quotePeriod.load(\"items\");\n  await context.sync();\n  if (quotePeriod.items.length > 0) {\n    quotePeriod.items[0].insertText(\"\u201d.\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// --- 3: remove the empty paragraph + the \"B\u00fasqueda de productos\" bullet ---\nparagraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet searchBulletIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"B\u00fasqueda de productos\") !== -1) {\n    searchBulletIndex = i;\n    break;\n  }\n}\n\nif (searchBulletIndex > 0) {\n  const precedingParagraph = paragraphs.items[searchBulletIndex - 1];\n  precedingParagraph.load(\"text\");\n  await context.sync();\n\n  if (precedingParagraph.text === \"\") {\n    // Delete the blank paragraph first (so the \"B\u00fasqueda\u2026\" paragraph's\n    // index doesn't shift before we grab it fresh), then reload and\n    // delete the \"B\u00fasqueda\u2026\" bullet itself.\n    precedingParagraph.delete();\n    await context.sync();\n\n    const refreshed = body.paragraphs;\n    refreshed.load(\"items/text\");\n    await context.sync();\n\n    for (let i = 0; i < refreshed.items.length; i++) {\n      if (refreshed.items[i].text.indexOf(\"B\u00fasqueda de productos\") !== -1) {\n        refreshed.items[i].delete();\n        break;\n      }\n    }\n    await context.sync();\n  } else {\n    // Fallback: preceding paragraph wasn't blank; just delete the bullet.\n    paragraphs.items[searchBulletIndex].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Fixed Search Box and updates Resumen Ejecutivo\n#\n# 1) Merge the split runs around the \"<quote> o <quote>\" text in the\n#    \"Edit Product\" / \"Delete Product\" bullet so that text lives in a\n#    single run.\n# 2) Merge the split runs around the trailing \"<quote>.\" at the end of\n#    that same bullet.\n# 3) Remove the empty paragraph and the \"Busqueda de productos...\" bullet\n#    at the end of the document (the search-box filter bullet was\n#    dropped from the write-up).\n\n$d = $word.ActiveDocument\n\n# --- 1: collapse the run-split around the closing/opening smart quotes\n#        (\"... Product\" o \"Delete ...\") ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"\u201d o \u201c\"\n$find1.Replacement.Text = \"\u201d o \u201c\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# --- 2: collapse the run-split around the trailing smart-quote + period\n#        (\"...Delete Product\".) ---\n# Scope the Find to the bullet paragraph itself so we only touch the\n# trailing quote+period pair, not every '\".' in the document.\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Edit Product*\" -and $t -like \"*Delete Product*\") {\n        $targetParagraph = $p\n    }\n}\n\nif ($targetParagraph -ne $null) {\n    $find2 = $targetParagraph.Range.Find\n    $find2.ClearFormatting()\n    $find2.Replacement.ClearFormatting()\n    $find2.Text = \"\u201d.\"\n    $find2.Replacement.Text = \"\u201d.\"\n    $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n}\n\n# --- 3: remove the empty paragraph + the \"Busqueda de productos\" bullet ---\n$searchBullet = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Se filtra el numero de productos*\") {\n        $searchBullet = $p\n    }\n}\n\nif ($searchBullet -ne $null) {\n    $precedingParagraph = $searchBullet.Previous()\n    # Paragraph.Range.Text includes the trailing paragraph mark, so an\n    # \"empty\" paragraph's text is \"`r\", not \"\".\n    if ($precedingParagraph.Range.Text.Trim() -eq \"\") {\n        # Delete the blank paragraph first, then re-fetch a fresh\n        # reference to the search bullet (its old COM object/range can\n        # go stale once a sibling paragraph is removed) and delete it.\n        $precedingParagraph.Range.Delete()\n\n        $searchBullet2 = $null\n        foreach ($p in $d.Paragraphs) {\n            if ($p.Range.Text -like \"*Se filtra el numero de productos*\") {\n                $searchBullet2 = $p\n            }\n        }\n        if ($searchBullet2 -ne $null) {\n            $searchBullet2.Range.Delete()\n        }\n    } else {\n        $searchBullet.Range.Delete()\n    }\n}\n"}
